$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2-34
# from serial date 45204 (2023-10-05) to 45207 (2023-10-08)
$ws.Range("C2:C34").Value = 45207
